# Apply "Discussions" feature requirement updates to the Functional Requirements
# workbook: mark the three existing "Discussions" rows (29-31, i.e. sheet rows
# 45-47) as TESTED (green highlight), and append four newly-tested Discussions
# requirements (REQ 161-164) at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Light-green "TESTED" fill color used throughout the sheet (RGB FF99FF99,
# stored by Excel as the BGR integer 10092441) together with the text number
# format ("@") that the REQ_ID / DEPENDENCY columns use.
$testedColor = 10092441

function Set-TestedRow {
    param(
        [int]$Row,
        [string]$ReqId,
        [string]$Category,
        [string]$Dependency,
        [string]$Description
    )

    if ($ReqId -ne $null) { $ws.Cells.Item($Row, 2).Value = $ReqId }
    if ($Category -ne $null) { $ws.Cells.Item($Row, 3).Value = $Category }
    if ($Dependency -ne $null) { $ws.Cells.Item($Row, 4).Value = $Dependency }
    if ($Description -ne $null) { $ws.Cells.Item($Row, 5).Value = $Description }

    # Column A / B / D carry the text number format, column C / E do not --
    # matching the existing "TESTED" rows elsewhere in the sheet.
    $ws.Cells.Item($Row, 1).Interior.Color = $testedColor
    $ws.Cells.Item($Row, 1).NumberFormat = "@"

    $ws.Cells.Item($Row, 2).Interior.Color = $testedColor
    $ws.Cells.Item($Row, 2).NumberFormat = "@"

    $ws.Cells.Item($Row, 3).Interior.Color = $testedColor

    $ws.Cells.Item($Row, 4).Interior.Color = $testedColor
    $ws.Cells.Item($Row, 4).NumberFormat = "@"

    $ws.Cells.Item($Row, 5).Interior.Color = $testedColor
}

# --- Mark the existing Discussions requirements (rows 45-47) as TESTED ------
# (NOTE: positional arguments are used throughout -- named arguments to
#  functions are not reliably bound by this host's PowerShell parser.)
Set-TestedRow 45 $null $null $null $null
Set-TestedRow 46 $null $null $null $null
Set-TestedRow 47 $null $null $null $null

# --- Append the newly tested Discussions requirements (rows 177-180) -------
Set-TestedRow 177 "161" "Discussions" $null "The post a reply button should redirect user to Reply creation page."
Set-TestedRow 178 "162" "Discussions" $null "Navigating to ./discussions/create/{mid} should take user to discussion creation page."
Set-TestedRow 179 "163" "Discussions" $null "Discussion creation form should be properly saved to the database."
Set-TestedRow 180 "164" "Discussions" $null "Most recent discussions should appear on landing page."

# --- Restore the on-screen selection to the new bottom of the table --------
$ws.Range("E184").Select() | Out-Null
